# Auto update: 2025-12-05 10:51:57
# Rotates the company/ticker labels in B2:C5 and refreshes the scoring
# metrics in D:K and N for rows 2-6 of the daily 방산(defense) screen.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: HYUNDAI ROTEM / 064350.KS ---
$ws.Range("B2").Value = "HYUNDAI ROTEM"
$ws.Range("C2").Value = "064350.KS"
$ws.Range("D2").Value = 180200
$ws.Range("E2").Value = 34.4
$ws.Range("F2").Value = 2.56
$ws.Range("G2").Value = 30
$ws.Range("H2").Value = 60
$ws.Range("I2").Value = 66
$ws.Range("J2").Value = 60
$ws.Range("K2").Value = 51.9
$ws.Range("N2").Value = 54.85170003294819

# --- Row 3: HANWHA AEROSPACE / 012450.KS ---
$ws.Range("B3").Value = "HANWHA AEROSPACE"
$ws.Range("C3").Value = "012450.KS"
$ws.Range("D3").Value = 864000
$ws.Range("E3").Value = 30.4
$ws.Range("F3").Value = 1.41
$ws.Range("G3").Value = 20
$ws.Range("H3").Value = 63
$ws.Range("I3").Value = 73
$ws.Range("J3").Value = 63
$ws.Range("K3").Value = 51.7
$ws.Range("N3").Value = 54.85170003294819

# --- Row 4: HANWHA SYSTEMS / 272210.KS ---
$ws.Range("B4").Value = "HANWHA SYSTEMS"
$ws.Range("C4").Value = "272210.KS"
$ws.Range("D4").Value = 46600
$ws.Range("E4").Value = 21.9
$ws.Range("F4").Value = 0.87
$ws.Range("G4").Value = 20
$ws.Range("H4").Value = 60
$ws.Range("I4").Value = 60
$ws.Range("J4").Value = 43
$ws.Range("K4").Value = 46.5
$ws.Range("N4").Value = 54.85170003294819

# --- Row 5: KOREA AEROSPACE / 047810.KS ---
$ws.Range("B5").Value = "KOREA AEROSPACE"
$ws.Range("C5").Value = "047810.KS"
$ws.Range("D5").Value = 104700
$ws.Range("E5").Value = 37.8
$ws.Range("F5").Value = -3.86
$ws.Range("G5").Value = 10
$ws.Range("H5").Value = 40
$ws.Range("I5").Value = 60
$ws.Range("J5").Value = 63
$ws.Range("K5").Value = 43.5
$ws.Range("N5").Value = 54.85170003294819

# --- Row 6: LIG Nex1 / 079550.KS (identity unchanged, metrics refreshed) ---
$ws.Range("D6").Value = 367000
$ws.Range("E6").Value = 26.8
$ws.Range("F6").Value = -4.18
$ws.Range("G6").Value = 10
$ws.Range("H6").Value = 56
$ws.Range("I6").Value = 46
$ws.Range("J6").Value = 46
$ws.Range("K6").Value = 37.9
$ws.Range("N6").Value = 54.85170003294819
